$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.980.12'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.916.37'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = '''324.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = '''0.4597'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.3826'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '''0.07702'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("D10").Value = '''0.9807'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").Value = '''22.22'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = '1.911.35'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '''5.688'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").Value = '''6.960'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").Value = '''0.07004'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '''84.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("D18").Value = '''0.000009479'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("D19").Value = '''16.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '28.982.45'
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("D22").Value = '''5.323'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.20%  '
$ws.Range("D23").Value = '''10.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.49%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.177.08'
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''2.085'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''158.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''19.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '''5.690'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '''117.82'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '''1.868'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.85%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.09310'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''0.8653'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.111'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.250'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''3.040'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.05704'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.158'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '''1.001'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02042'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''3.054'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.19%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''7.504'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.5512'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.1752'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''9.382'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '''0.000002832'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.82%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''2.180'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.32%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.5191'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''11.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.06899'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.783'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '''110.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
